$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 299.66666
$ws.Cells.Item(18, 9).Value = 299.66666
$ws.Cells.Item(18, 11).Value = 299.66666
$ws.Cells.Item(18, 13).Value = -15.66665999999998
$ws.Cells.Item(51, 8).Value = 9999.5
$ws.Cells.Item(98, 8).Value = 2500
$ws.Cells.Item(98, 9).Value = 2500
$ws.Cells.Item(98, 11).Value = 2500
$ws.Cells.Item(98, 13).Value = -1002
$ws.Cells.Item(100, 8).Value = 746.5
$ws.Cells.Item(100, 9).Value = 746.5
$ws.Cells.Item(100, 11).Value = 746.5
$ws.Cells.Item(100, 13).Value = -205.5
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 13).ClearContents()
$ws.Cells.Item(112, 8).Value = 1879
$ws.Cells.Item(112, 10).Value = 2049.7646
$ws.Cells.Item(112, 12).Value = 6149.293799999999
$ws.Cells.Item(112, 14).Value = -8365.293799999999
$ws.Cells.Item(122, 8).Value = 2500
$ws.Cells.Item(122, 9).Value = 2500
$ws.Cells.Item(122, 11).Value = 7500
$ws.Cells.Item(122, 13).Value = -5050
$ws.Cells.Item(138, 8).Value = 4528.477
$ws.Cells.Item(138, 10).Value = 4818.914
$ws.Cells.Item(138, 12).Value = 14456.742
$ws.Cells.Item(138, 14).Value = -24736.742
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 1326.909
$ws.Cells.Item(5, 9).Value = 955.875
$ws.Cells.Item(5, 10).Value = 2316.3333
$ws.Cells.Item(5, 11).Value = 955.875
$ws.Cells.Item(5, 12).Value = 2316.3333
$ws.Cells.Item(5, 13).Value = -843.875
$ws.Cells.Item(5, 14).Value = -2540.3333
$ws.Cells.Item(33, 8).Value = 3506.25
$ws.Cells.Item(33, 9).Value = 1341.6666
$ws.Cells.Item(33, 10).Value = 10000
$ws.Cells.Item(33, 11).Value = 1341.6666
$ws.Cells.Item(33, 12).Value = 10000
$ws.Cells.Item(33, 13).Value = -1012.6666
$ws.Cells.Item(33, 14).Value = -10658
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(86, 8).Value = 37213.5
$ws.Cells.Item(86, 9).Value = 37313
$ws.Cells.Item(86, 10).Value = 37114
$ws.Cells.Item(86, 11).Value = 37313
$ws.Cells.Item(86, 12).Value = 37114
$ws.Cells.Item(86, 13).Value = -36127
$ws.Cells.Item(86, 14).Value = -39486
$ws.Cells.Item(89, 8).Value = 37213.5
$ws.Cells.Item(89, 9).Value = 37313
$ws.Cells.Item(89, 10).Value = 37114
$ws.Cells.Item(89, 11).Value = 111939
$ws.Cells.Item(89, 12).Value = 111342
$ws.Cells.Item(89, 13).Value = -106011
$ws.Cells.Item(89, 14).Value = -123198
$ws.Cells.Item(132, 8).Value = 2708.2144
$ws.Cells.Item(132, 9).Value = 772.4286
$ws.Cells.Item(132, 10).Value = 4644
$ws.Cells.Item(132, 11).Value = 2317.2858
$ws.Cells.Item(132, 12).Value = 13932
$ws.Cells.Item(132, 13).Value = 212.7142000000003
$ws.Cells.Item(132, 14).Value = -18992
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 1326.909
$ws.Cells.Item(4, 9).Value = 955.875
$ws.Cells.Item(4, 10).Value = 2316.3333
$ws.Cells.Item(4, 11).Value = 955.875
$ws.Cells.Item(4, 12).Value = 2316.3333
$ws.Cells.Item(4, 13).Value = -840.875
$ws.Cells.Item(4, 14).Value = -2546.3333
$ws.Cells.Item(22, 8).Value = 284
$ws.Cells.Item(22, 9).Value = 292.66666
$ws.Cells.Item(22, 11).Value = 292.66666
$ws.Cells.Item(22, 13).Value = -119.66666
$ws.Cells.Item(80, 8).Value = 1127.6
$ws.Cells.Item(80, 9).Value = 92.36364
$ws.Cells.Item(80, 10).Value = 3974.5
$ws.Cells.Item(80, 11).Value = 92.36364
$ws.Cells.Item(80, 12).Value = 3974.5
$ws.Cells.Item(80, 13).Value = 905.63636
$ws.Cells.Item(80, 14).Value = -5970.5
$ws.Cells.Item(83, 8).Value = 1127.6
$ws.Cells.Item(83, 9).Value = 92.36364
$ws.Cells.Item(83, 10).Value = 3974.5
$ws.Cells.Item(83, 11).Value = 461.8182
$ws.Cells.Item(83, 12).Value = 19872.5
$ws.Cells.Item(83, 13).Value = 4530.1818
$ws.Cells.Item(83, 14).Value = -29856.5
$ws.Cells.Item(86, 8).Value = 4779.4
$ws.Cells.Item(89, 8).Value = 4779.4
$ws.Cells.Item(105, 8).Value = 4369.2
$ws.Cells.Item(105, 9).Value = 4505
$ws.Cells.Item(105, 10).Value = 4335.25
$ws.Cells.Item(105, 11).Value = 4505
$ws.Cells.Item(105, 12).Value = 4335.25
$ws.Cells.Item(105, 13).Value = -2758
$ws.Cells.Item(105, 14).Value = -7829.25
$ws.Cells.Item(112, 8).Value = 24466.334
$ws.Cells.Item(112, 10).Value = 24466.334
$ws.Cells.Item(112, 12).Value = 24466.334
$ws.Cells.Item(112, 14).Value = -27420.334
$ws.Cells.Item(134, 8).Value = 4890.6665
$ws.Cells.Item(134, 10).Value = 4336
$ws.Cells.Item(134, 12).Value = 13008
$ws.Cells.Item(134, 14).Value = -18078
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(8, 8).Value = 1749
$ws.Cells.Item(8, 9).Value = 1749
$ws.Cells.Item(8, 11).Value = 1749
$ws.Cells.Item(8, 13).Value = -1609
$ws.Cells.Item(22, 8).Value = 725.3333
$ws.Cells.Item(22, 9).Value = 595
$ws.Cells.Item(22, 11).Value = 595
$ws.Cells.Item(22, 13).Value = -245
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 13).ClearContents()
$ws.Cells.Item(31, 8).Value = 4449.4
$ws.Cells.Item(31, 9).Value = 3057.3333
$ws.Cells.Item(31, 10).Value = 6537.5
$ws.Cells.Item(31, 11).Value = 3057.3333
$ws.Cells.Item(31, 12).Value = 6537.5
$ws.Cells.Item(31, 13).Value = -2762.3333
$ws.Cells.Item(31, 14).Value = -7127.5
$ws.Cells.Item(34, 8).Value = 4449.4
$ws.Cells.Item(34, 9).Value = 3057.3333
$ws.Cells.Item(34, 10).Value = 6537.5
$ws.Cells.Item(34, 11).Value = 3057.3333
$ws.Cells.Item(34, 12).Value = 6537.5
$ws.Cells.Item(34, 13).Value = -2855.3333
$ws.Cells.Item(34, 14).Value = -6941.5
$ws.Cells.Item(135, 8).Value = 72316.164
$ws.Cells.Item(135, 10).Value = 72316.164
$ws.Cells.Item(135, 12).Value = 72316.164
$ws.Cells.Item(135, 14).Value = -82456.164
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 167.5
$ws.Cells.Item(10, 9).Value = 167.5
$ws.Cells.Item(10, 11).Value = 502.5
$ws.Cells.Item(10, 13).Value = -363.5
$ws.Cells.Item(45, 8).Value = 5479.7144
$ws.Cells.Item(45, 10).Value = 5279.6
$ws.Cells.Item(45, 12).Value = 15838.8
$ws.Cells.Item(45, 14).Value = -16902.8
$ws.Cells.Item(50, 8).Value = 3691
$ws.Cells.Item(50, 9).Value = 5074.25
$ws.Cells.Item(50, 11).Value = 15222.75
$ws.Cells.Item(50, 13).Value = -14741.75
$ws.Cells.Item(53, 8).Value = 3691
$ws.Cells.Item(53, 9).Value = 5074.25
$ws.Cells.Item(53, 11).Value = 15222.75
$ws.Cells.Item(53, 13).Value = -14741.75
$ws.Cells.Item(81, 8).Value = 56481.25
$ws.Cells.Item(81, 9).Value = 38637.668
$ws.Cells.Item(81, 10).Value = 110012
$ws.Cells.Item(81, 11).Value = 115913.004
$ws.Cells.Item(81, 12).Value = 330036
$ws.Cells.Item(81, 13).Value = -114790.004
$ws.Cells.Item(81, 14).Value = -332282
$ws.Cells.Item(84, 8).Value = 56481.25
$ws.Cells.Item(84, 9).Value = 38637.668
$ws.Cells.Item(84, 10).Value = 110012
$ws.Cells.Item(84, 11).Value = 347739.012
$ws.Cells.Item(84, 12).Value = 990108
$ws.Cells.Item(84, 13).Value = -342123.012
$ws.Cells.Item(84, 14).Value = -1001340
$ws.Cells.Item(122, 8).Value = 816.25
$ws.Cells.Item(122, 9).Value = 1060.5
$ws.Cells.Item(122, 10).Value = 572
$ws.Cells.Item(122, 11).Value = 9544.5
$ws.Cells.Item(122, 12).Value = 5148
$ws.Cells.Item(122, 13).Value = -7094.5
$ws.Cells.Item(122, 14).Value = -10048
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 10239.28
$ws.Cells.Item(80, 9).Value = 4299.4614
$ws.Cells.Item(80, 11).Value = 4299.4614
$ws.Cells.Item(80, 13).Value = -3301.4614
$ws.Cells.Item(83, 8).Value = 10239.28
$ws.Cells.Item(83, 9).Value = 4299.4614
$ws.Cells.Item(83, 11).Value = 21497.307
$ws.Cells.Item(83, 13).Value = -16505.307
$ws.Cells.Item(122, 8).Value = 203157.2
$ws.Cells.Item(122, 9).Value = 336665.66
$ws.Cells.Item(122, 10).Value = 2894.5
$ws.Cells.Item(122, 11).Value = 1009996.98
$ws.Cells.Item(122, 12).Value = 8683.5
$ws.Cells.Item(122, 13).Value = -1007546.98
$ws.Cells.Item(122, 14).Value = -13583.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 5000
$ws.Cells.Item(13, 9).Value = 5000
$ws.Cells.Item(13, 11).Value = 5000
$ws.Cells.Item(13, 13).Value = -4860
$ws.Cells.Item(25, 8).Value = 8000
$ws.Cells.Item(25, 9).Value = 20000
$ws.Cells.Item(25, 11).Value = 20000
$ws.Cells.Item(25, 13).Value = -19770
$ws.Cells.Item(34, 8).Value = 10830
$ws.Cells.Item(34, 10).Value = 10830
$ws.Cells.Item(34, 12).Value = 10830
$ws.Cells.Item(34, 14).Value = -11174
$ws.Cells.Item(46, 8).Value = 1763
$ws.Cells.Item(46, 9).Value = 2056.125
$ws.Cells.Item(46, 11).Value = 2056.125
$ws.Cells.Item(46, 13).Value = -1868.125
$ws.Cells.Item(55, 8).Value = 493.3846
$ws.Cells.Item(55, 9).Value = 421
$ws.Cells.Item(55, 11).Value = 421
$ws.Cells.Item(55, 13).Value = -248
$ws.Cells.Item(70, 8).Value = 48037.668
$ws.Cells.Item(70, 9).Value = 25950
$ws.Cells.Item(70, 11).Value = 25950
$ws.Cells.Item(70, 13).Value = -25680
$ws.Cells.Item(73, 8).Value = 48037.668
$ws.Cells.Item(73, 9).Value = 25950
$ws.Cells.Item(73, 11).Value = 25950
$ws.Cells.Item(73, 13).Value = -25014
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1494.3334
$ws.Cells.Item(81, 10).Value = 1494.5
$ws.Cells.Item(81, 12).Value = 2989
$ws.Cells.Item(81, 14).Value = -5111
$ws.Cells.Item(84, 8).Value = 1494.3334
$ws.Cells.Item(84, 10).Value = 1494.5
$ws.Cells.Item(84, 12).Value = 14945
$ws.Cells.Item(84, 14).Value = -25553
$ws.Cells.Item(107, 8).Value = 922
$ws.Cells.Item(107, 9).Value = 389
$ws.Cells.Item(107, 10).Value = 1759.5714
$ws.Cells.Item(107, 11).Value = 1167
$ws.Cells.Item(107, 12).Value = 5278.7142
$ws.Cells.Item(107, 13).Value = 753
$ws.Cells.Item(107, 14).Value = -9118.7142
$ws.Cells.Item(113, 8).Value = 1258.5555
$ws.Cells.Item(113, 10).Value = 1749.5
$ws.Cells.Item(113, 12).Value = 5248.5
$ws.Cells.Item(113, 14).Value = -9588.5
$ws.Cells.Item(122, 8).Value = 5039.028
$ws.Cells.Item(122, 9).Value = 5113.607
$ws.Cells.Item(122, 10).Value = 4778
$ws.Cells.Item(122, 11).Value = 15340.821
$ws.Cells.Item(122, 12).Value = 14334
$ws.Cells.Item(122, 13).Value = -12890.821
$ws.Cells.Item(122, 14).Value = -19234
$ws.Cells.Item(125, 8).Value = 53332.668
$ws.Cells.Item(125, 10).Value = 53332.668
$ws.Cells.Item(125, 12).Value = 53332.668
$ws.Cells.Item(125, 14).Value = -63172.668
$ws.Cells.Item(136, 8).Value = 4122.25
$ws.Cells.Item(136, 9).Value = 4141.143
$ws.Cells.Item(136, 11).Value = 12423.429
$ws.Cells.Item(136, 13).Value = -9873.429
